$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Population_name" values in B2 and B4 to drop the trailing
# date stamp, matching the now-current value used elsewhere (column C).
$ws.Range("B2").Value = "NewImportLogic_3 - Test_Automation_3"
$ws.Range("B4").Value = "NewImportLogic_3 - Test_Automation_3"

# Reset the view: scroll back to top-left and move the selection to B5.
$ws.Range("A1").Select()
$ws.Range("B5").Select()
